# Insert a new weekly price record for "Piña" (Vega Monumental Concepción)
# at row 241, pushing the existing rows 241:274 down to 242:275.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 241..274 down one row (Excel copies formatting from the row
# above, which is how the D column keeps its date style).
$ws.Rows("241").Insert()

# Populate the newly inserted row 241 with the new record's data.
$ws.Range("A241").Value = 11
$ws.Range("B241").Value = 'Vega Monumental Concepción'
$ws.Range("C241").Value = 'Bíobío'
$ws.Range("D241").Value = 45077
$ws.Range("E241").Value = 8
$ws.Range("F241").Value = 'Fruta'
$ws.Range("G241").Value = 100108
$ws.Range("H241").Value = 'Tropicales y subtropicales'
$ws.Range("I241").Value = 100108005
$ws.Range("J241").Value = 'Piña'
$ws.Range("K241").Value = 'Sin especificar'
$ws.Range("L241").Value = 'Segunda'
$ws.Range("M241").Value = 140
$ws.Range("N241").Value = 14000
$ws.Range("O241").Value = 15000
$ws.Range("P241").Value = 14571
$ws.Range("Q241").Value = '$/caja 14 unidades'
$ws.Range("R241").Value = 'Ecuador'
$ws.Range("S241").Value = 1041
$ws.Range("T241").Value = 14
